$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title / source strings (January -> February 2017)
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("C32").Value = "Source: Short-Term Energy Outlook, February 2017."

# Row 27 - OECD* (D:G raw data; I:K are shared formulas E-D, F-E, G-F and recalc automatically)
$ws.Range("D27").Value = 46.327475290000002
$ws.Range("E27").Value = 46.631664399999998
$ws.Range("F27").Value = 47.026988566999997
$ws.Range("G27").Value = 47.257554519999999

# Row 28 - Non-OECD Asia
$ws.Range("D28").Value = 24.23996868
$ws.Range("E28").Value = 25.127751891999999
$ws.Range("F28").Value = 25.900298849000002
$ws.Range("G28").Value = 26.661100000000001

# Row 29 - FSU and Eastern Europe
$ws.Range("D29").Value = 5.56511083944
$ws.Range("E29").Value = 5.5656669179500007
$ws.Range("F29").Value = 5.5777945163
$ws.Range("G29").Value = 5.5914392981800001

# Row 31 - World (D30:G30 "Other" is a formula D31-D27-D28-D29 and recalcs automatically)
$ws.Range("D31").Value = 95.060221308999999
$ws.Range("E31").Value = 96.466253347000006
$ws.Range("F31").Value = 98.088077510000005
$ws.Range("G31").Value = 99.545306836999998
